$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Sheet1"

# 2. Fill in new category rows (B4 plus new rows 5-8)
$ws.Range("B4").Value = "Thể loại 1"
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Thể loại 2"
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Thể loại 3"
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Thể loại 4"
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Thể loại 5"

# 3. Propagate B4's format (font/border/alignment) down to B5:B8 before adding fill
$ws.Range("B4").Copy()
$ws.Range("B5:B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. White fill for column B data block + A4
$ws.Range("A4").Interior.Color = 16777215
$ws.Range("B4:B8").Interior.Color = 16777215

# 5. Thin border around the new data rows in column A
$rngBorder = $ws.Range("A5:A8")
$rngBorder.Borders.LineStyle = 1
$rngBorder.Borders.Weight = 2

# 6. Column B width
$ws.Columns("B").ColumnWidth = 32.285714285714285

# 7. Selection moves to E4
$ws.Range("E4").Select() | Out-Null
